$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)

# Insert 5 new rows before the existing row 5 (the "if" branch block),
# pushing the rest of the flow down to make room for the new "else" branch.
$ws.Rows("5:9").Insert()

# New "if" condition block (mirrors the existing one lower down)
$ws.Range("B5").Value = "if"
$ws.Range("C5").Value = "data('test3') > 100"

# Forgotten "note" row describing what happens when the condition is true
$ws.Range("D6").Value = "note"
$ws.Range("H6").Value = "Will exit section"

$ws.Range("B7").Value = "else"

# Forgotten "note" row describing what happens when the condition is false
$ws.Range("D8").Value = "note"
$ws.Range("H8").Value = "Will not exit section"

$ws.Range("B9").Value = "end if"

# Update the selection to reflect where the author was last working
$ws.Range("D9").Select()
